# Adds two new columns, I ("I0") and J ("IF"), to Sheet1, with header
# labels in row 1 and per-row numeric data in rows 2-58.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (copy formatting from the existing "IP" header cell so the
# new headers pick up the same bold/border/alignment style, s="1")
$ws.Cells.Item(1, 9).Value = "I0"
$ws.Cells.Item(1, 10).Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows: row index, I value, J value
$data = @(
    @(2, 1, 3),
    @(3, 6, 7),
    @(4, 8, 8),
    @(5, 6, 6),
    @(6, 8, 8),
    @(7, 8, 8),
    @(8, 4, 6),
    @(9, 9, 9),
    @(10, 8, 8),
    @(11, 9, 9),
    @(12, 8, 8),
    @(13, 6, 6),
    @(14, 8, 8),
    @(15, 7, 7),
    @(16, 6, 7),
    @(17, 7, 8),
    @(18, 7, 7),
    @(19, 8, 8),
    @(20, 9, 9),
    @(21, 8, 8),
    @(22, 9, 9),
    @(23, 8, 9),
    @(24, 9, 9),
    @(25, 9, 9),
    @(26, 8, 8),
    @(27, 8, 8),
    @(28, 8, 8),
    @(29, 7, 8),
    @(30, 9, 9),
    @(31, 6, 8),
    @(32, 5, 7),
    @(33, 8, 9),
    @(34, 10, 10),
    @(35, 6, 6),
    @(36, 8, 8),
    @(37, 9, 9),
    @(38, 8, 8),
    @(39, 6, 7),
    @(40, 7, 9),
    @(41, 6, 8),
    @(42, 7, 9),
    @(43, 6, 9),
    @(44, 1, 5),
    @(45, 1, 6),
    @(46, 4, 8),
    @(47, 1, 6),
    @(48, 1, 5),
    @(49, 5, 8),
    @(50, 1, 4),
    @(51, 1, 7),
    @(52, 1, 5),
    @(53, 1, 5),
    @(54, 1, 5),
    @(55, 1, 5),
    @(56, 1, 5),
    @(57, 1, 4),
    @(58, 1, 3)
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
